$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A (Fecha) timestamps for rows 2-10 to the new date (2025-08-07 series)
$ws.Range("A2").Value = 45876.00017361111
$ws.Range("A3").Value = 45876.04184027778
$ws.Range("A4").Value = 45876.08350694444
$ws.Range("A5").Value = 45876.12517361111
$ws.Range("A6").Value = 45876.16684027778
$ws.Range("A7").Value = 45876.20850694444
$ws.Range("A8").Value = 45876.25017361111
$ws.Range("A9").Value = 45876.29184027778
$ws.Range("A10").Value = 45876.33350694444

# Tiny correction to row 11's timestamp
$ws.Range("A11").Value = 45876.37517238426

# Add new row 12 with the latest reading
$ws.Range("A12").Value = 45876.41687603376
$ws.Range("B12").Value = 2025
$ws.Range("C12").Value = 28
$ws.Range("D12").Value = 16.06
$ws.Range("E12").Value = 93.40000000000001
$ws.Range("F12").Value = 472.85
$ws.Range("G12").Value = 4.23
$ws.Range("H12").Value = "SE"
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = "10:00:18"

# Match the date style/number format used in column A for the other rows
$ws.Range("A12").NumberFormat = $ws.Range("A11").NumberFormat
